# Apply SHGC revision (0.57 -> 0.26) and updated peak heat/cool results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update SHGC label value (Part 1 summary block) ---
$ws.Range("E1").Value = 0.26

# --- Update peak heating/cooling load + capacity results (first block, rows 4-7) ---
$ws.Range("B4").Value = 14.1
$ws.Range("C4").Value = 2.15
$ws.Range("E4").Value = 0.92

$ws.Range("B5").Value = 14.1
$ws.Range("C5").Value = 2.15
$ws.Range("E5").Value = 0.92

$ws.Range("B6").Value = 14.56
$ws.Range("C6").Value = 2.45
$ws.Range("E6").Value = 1.04

$ws.Range("B7").Value = 14.86
$ws.Range("C7").Value = 2.66
$ws.Range("E7").Value = 1.12

# --- Update duplicate results block (rows 16-19) ---
$ws.Range("B16").Value = 14.1
$ws.Range("C16").Value = 2.15
$ws.Range("E16").Value = 0.92

$ws.Range("B17").Value = 14.1
$ws.Range("C17").Value = 2.15
$ws.Range("E17").Value = 0.92

$ws.Range("B18").Value = 14.56
$ws.Range("C18").Value = 2.45
$ws.Range("E18").Value = 1.04

$ws.Range("B19").Value = 14.86
$ws.Range("C19").Value = 2.66
$ws.Range("E19").Value = 1.12

# --- Restore the sheet view to top-left (remove scrolled/selected state) ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()

$excel.CalculateFullRebuild()
